$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" value would otherwise be auto-recognised by Excel
# as a number (losing the original text formatting / trailing zeros). Force
# text storage for just these cells, then restore the default style so no
# visible formatting change is introduced.
$numericLookingCells = @(
"D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D38", "D39", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D49"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.276.86"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "1.830.82"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").Value = "235.23"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "0.6041"
$ws.Range("E6").Value = "  -3.80%  "

$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").Value = "0.07035"
$ws.Range("E8").Value = "  -5.50%  "

$ws.Range("D9").Value = "0.2794"
$ws.Range("E9").Value = "  -3.61%  "

$ws.Range("D10").Value = "23.49"
$ws.Range("E10").Value = "  -5.75%  "

$ws.Range("D11").Value = "0.07663"
$ws.Range("E11").Value = "  -0.66%  "

$ws.Range("D12").Value = "1.832.10"
$ws.Range("E12").Value = "  -0.86%  "

$ws.Range("D13").Value = "4.805"
$ws.Range("E13").Value = "  -3.39%  "

$ws.Range("D14").Value = "0.000009904"
$ws.Range("E14").Value = "  -3.46%  "

$ws.Range("D15").Value = "0.6260"
$ws.Range("E15").Value = "  -7.44%  "

$ws.Range("D16").Value = "79.22"
$ws.Range("E16").Value = "  -3.20%  "

$ws.Range("D17").Value = "29.251.54"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "5.842"
$ws.Range("E18").Value = "  -6.25%  "

$ws.Range("D19").Value = "224.98"
$ws.Range("E19").Value = "  -3.48%  "

$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("E21").Value = "  -5.00%  "

$ws.Range("D22").Value = "7.011"
$ws.Range("E22").Value = "  -4.21%  "

$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "156.68"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D25").Value = "8.006"
$ws.Range("E25").Value = "  -5.70%  "

$ws.Range("D26").Value = "0.1297"
$ws.Range("E26").Value = "  -4.08%  "

$ws.Range("D27").Value = "16.55"
$ws.Range("E27").Value = "  -4.65%  "

$ws.Range("D28").Value = "1.479"
$ws.Range("E28").Value = "  +0.94%  "

$ws.Range("D29").Value = "0.06185"
$ws.Range("E29").Value = "  -12.90%  "

$ws.Range("E30").Value = "  -2.47%  "

$ws.Range("D31").Value = "3.829"
$ws.Range("E31").Value = "  -5.20%  "

$ws.Range("D32").Value = "3.799"
$ws.Range("E32").Value = "  -6.14%  "

$ws.Range("D33").Value = "1.122"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("D35").Value = "0.6453"
$ws.Range("E35").Value = "  -7.61%  "

$ws.Range("D36").Value = "2.546"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("D37").Value = "1.223.64"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("D38").Value = "2.740"
$ws.Range("E38").Value = "  -2.66%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "6.559"
$ws.Range("E39").Value = "  -6.28%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01735"
$ws.Range("E40").Value = "  -5.58%  "

$ws.Range("D41").Value = "0.8986"
$ws.Range("E41").Value = "  -6.12%  "

$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").Value = "1.988.11"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("D44").Value = "100.66"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").Value = "62.57"
$ws.Range("E45").Value = "  -4.47%  "

$ws.Range("D46").Value = "0.00000000116"
$ws.Range("E46").Value = "  -3.02%  "

$ws.Range("D47").Value = "8.537"
$ws.Range("E47").Value = "  -4.58%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.581"
$ws.Range("E48").Value = "  -8.52%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.4563"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("E50").Value = "  -2.58%  "

$ws.Range("E51").Value = "  -7.76%  "

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
